$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "file" column (A) holds per-country source-file ids
# (e.g. 87811004_1121_AU). Sourcing is being made "universal": the data
# rows (2-23) are re-sorted into ascending alphabetical order by that
# column, keeping every row's B:E values (r_count/currency/sum/
# built_in_total) intact together with their original cell formatting.
$dataRange = $ws.Range("A2:E23")
$keyRange = $ws.Range("A2:A23")

$dataRange.Sort($keyRange)
